$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed data values in row 3
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Reflect the new active cell / selection (E3) as recorded in the saved view state
$ws.Range("E3").Select()
